$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Burndown update: one more hour was logged against "Gestão de git e docs
# (SM)" on Dia 4 (column E, row 24). The "Restante" row (25) recalculates
# automatically via its SUM-based formulas.
$ws.Range("E24").Value = 2

# Move the sheet's active selection to E23, matching where the author left
# the cursor after the edit.
$ws.Activate()
$ws.Range("E23").Select()
